$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Remove-HyperlinkAt($sheet, $addr) {
    $found = $null
    foreach ($h in $sheet.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $found = $h
        }
    }
    if ($found -ne $null) {
        $found.Delete()
    }
}

# ---- Header row (row 1) ----
$ws.Range("D1").Value = "file_1"
$ws.Range("F1").Value = "subject"

# ---- Row 2: new "subject" column value ----
$ws.Range("F2").Value = "Software developer internship"

# ---- Row 3: new "subject" column value ----
$ws.Range("F3").Value = "Software engineering internship"

# ---- Row 4: template file switched back to the docx, new "subject" value ----
$ws.Range("D4").Value = "test_template.docx"
$ws.Range("F4").Value = "data science internship"

# ---- Row 5: drop the recruiter/email hyperlinks + their cell text ----
Remove-HyperlinkAt $ws '$C$5'
Remove-HyperlinkAt $ws '$E$5'
$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()

# ---- Row 6: brand new row ----
$ws.Range("A6").Value = "Company E"
$ws.Range("B6").Value = "Something else"

$ws.Range("C6").Value = "some_other_email@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:some_other_email@gmail.com")
$ws.Range("C6").Style = "Hyperlink"

$ws.Range("E6").Value = "rayanakhtar120330@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:rayanakhtar120330@gmail.com")
$ws.Range("E6").Style = "Hyperlink"

# ---- Column widths (C widened, F newly sized) ----
$ws.Columns("C").ColumnWidth = 26.191616766467067
$ws.Columns("F").ColumnWidth = 30.583333333954215

# ---- Selection moved to F5 ----
$ws.Range("F5").Select()
